$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ngf"
$ws.Range("C2").Value = "Ngfr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7246536666666668
$ws.Range("H2").Value = 2.173961
$ws.Range("I2").Value = 0.1791272621505297
$ws.Range("J2").Value = 0.1791272621505298
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6946430000000001
$ws.Range("N2").Value = 2.083929
$ws.Range("O2").Value = 0.1269399741689062
$ws.Range("P2").Value = 0.1269399741689062
$ws.Range("Q2").Value = 0.5033755969743334
$ws.Range("R2").Value = 4.530380372769002
$ws.Range("S2").Value = 0.02273841003033514
$ws.Range("T2").Value = 0.02273841003033514

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ngf"
$ws.Range("C3").Value = "Ngfr"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7246536666666668
$ws.Range("H3").Value = 2.173961
$ws.Range("I3").Value = 0.1791272621505297
$ws.Range("J3").Value = 0.1791272621505298
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05042666666666667
$ws.Range("N3").Value = 0.15128
$ws.Range("O3").Value = 0.009215035297398391
$ws.Range("P3").Value = 0.009215035297398391
$ws.Range("Q3").Value = 0.03654186889777778
$ws.Range("R3").Value = 0.32887682008
$ws.Range("S3").Value = 0.001650664043443466
$ws.Range("T3").Value = 0.001650664043443467

# Row 4
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Ngf"
$ws.Range("C4").Value = "Ngfr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7246536666666668
$ws.Range("H4").Value = 2.173961
$ws.Range("I4").Value = 0.1791272621505297
$ws.Range("J4").Value = 0.1791272621505298
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.727146666666667
$ws.Range("N4").Value = 14.18144
$ws.Range("O4").Value = 0.8638449905336953
$ws.Range("P4").Value = 0.8638449905336955
$ws.Range("Q4").Value = 3.425544164871112
$ws.Range("R4").Value = 30.82989748384
$ws.Range("S4").Value = 0.1547381880767511
$ws.Range("T4").Value = 0.1547381880767512

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Ngf"
$ws.Range("C5").Value = "Ngfr"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.320814666666667
$ws.Range("H5").Value = 9.962444
$ws.Range("I5").Value = 0.8208727378494701
$ws.Range("J5").Value = 0.8208727378494702
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6946430000000001
$ws.Range("N5").Value = 2.083929
$ws.Range("O5").Value = 0.1269399741689062
$ws.Range("P5").Value = 0.1269399741689062
$ws.Range("Q5").Value = 2.306780662497334
$ws.Range("R5").Value = 20.761025962476
$ws.Range("S5").Value = 0.1042015641385711
$ws.Range("T5").Value = 0.1042015641385711

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Ngf"
$ws.Range("C6").Value = "Ngfr"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.320814666666667
$ws.Range("H6").Value = 9.962444
$ws.Range("I6").Value = 0.8208727378494701
$ws.Range("J6").Value = 0.8208727378494702
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05042666666666667
$ws.Range("N6").Value = 0.15128
$ws.Range("O6").Value = 0.009215035297398391
$ws.Range("P6").Value = 0.009215035297398391
$ws.Range("Q6").Value = 0.1674576142577778
$ws.Range("R6").Value = 1.50711852832
$ws.Range("S6").Value = 0.007564371253954923
$ws.Range("T6").Value = 0.007564371253954924

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Ngf"
$ws.Range("C7").Value = "Ngfr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.320814666666667
$ws.Range("H7").Value = 9.962444
$ws.Range("I7").Value = 0.8208727378494701
$ws.Range("J7").Value = 0.8208727378494702
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.727146666666667
$ws.Range("N7").Value = 14.18144
$ws.Range("O7").Value = 0.8638449905336953
$ws.Range("P7").Value = 0.8638449905336955
$ws.Range("Q7").Value = 15.69797798215111
$ws.Range("R7").Value = 141.28180183936
$ws.Range("S7").Value = 0.7091068024569441
$ws.Range("T7").Value = 0.7091068024569442
